$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 16: fix the text of the ratio metric names (word order swapped).
$ws.Range("A16").Value = "ratio_last_first, ratio_intermediate_first, ratio_multiple_single"

# Row 17: remove the special "Lucida Console" font formatting from A17 (copy
# the plain formatting used elsewhere in column A), and give the row an
# explicit custom height.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A17").RowHeight = 36

# Selection moved from B17 to A17.
$ws.Range("A17").Select() | Out-Null

$wb.Save()
